$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TeamStats")

# Add the missing AVERAGE formula for the velocity column (B16), matching
# the SUM/AVERAGE formulas already present for the other columns.
$ws.Range("B16").Formula = "=AVERAGE(B2:B14)"

# Update the active selection to match the latest edit location.
$ws.Range("B16:E16").Select()

$wb.Save()
